# Applies the monthly "Actualización automática" update to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 3 (ARCOS GOMEZ CONSTRUCCIONES CIA. LTDA.): clear PIEDRA SINTERIZADA / PORCELANATO values
$ws1.Range("L3").Value = 0
$ws1.Range("M3").Value = 0

# Row 12 totals/counters reflecting the row-3 change
$ws1.Range("L12").Value = "0 de 10"
$ws1.Range("M12").Value = "0 de 10"

# --- Sheet "VENTA MENSUAL" -------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths adjust slightly as the month columns roll forward
# (subtract the fixed COM padding offset so the saved OOXML width is exactly 10 / 12)
$colWidthPad = 0.8333333333333339
$ws2.Columns.Item(3).ColumnWidth = 10 - $colWidthPad
$ws2.Columns.Item(4).ColumnWidth = 12 - $colWidthPad

# Header: months roll forward by one (abril/mayo/junio/julio -> mayo/junio/julio/agosto)
$ws2.Range("C1").Value = "mayo"
$ws2.Range("D1").Value = "junio"
$ws2.Range("E1").Value = "julio"
$ws2.Range("F1").Value = "agosto"

# Row 3 - ARCOS GOMEZ CONSTRUCCIONES CIA. LTDA.
$ws2.Range("D3").Value = 832
$ws2.Range("E3").Value = 594.47
$ws2.Range("F3").Value = 0

# Row 5 - CARRION ALVAREZ MARIO ANDRES
$ws2.Range("D5").Value = 155.38
$ws2.Range("E5").Value = 0

# Row 10 - VACA PANCHI DORYS CAROLINA
$ws2.Range("C10").Value = 0
$ws2.Range("D10").Value = 10.44
$ws2.Range("E10").Value = 0

# Row 12 - totals
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 997.8200000000001
$ws2.Range("E12").Value = 594.47
$ws2.Range("F12").Value = 0

$wb.Save()
